# Update the order date (A1) and recompute the two price rows (D33, D34)
# on the active sheet ("Hoja1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436

$ws.Range("D33").Value = 64.13500000000001
$ws.Range("D34").Value = 50.407
